# Auto-generated Excel COM-interop script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item('Sheet1')

$ws.Range('D2').Value = '''69.650.71'
$ws.Range('E2').Value = '  -1.62%  '
$ws.Range('D3').Value = '''3.504.78'
$ws.Range('E3').Value = '  -1.40%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '''616.30'
$ws.Range('E5').Value = '  +5.85%  '
$ws.Range('D6').Value = '''191.38'
$ws.Range('E6').Value = '  +2.24%  '
$ws.Range('D7').Value = '''0.627'
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('E9').Value = '  -2.75%  '
$ws.Range('D10').Value = '''0.654'
$ws.Range('E10').Value = '  +0.79%  '
$ws.Range('D11').Value = '''53.18'
$ws.Range('E11').Value = '  -2.07%  '
$ws.Range('D12').Value = '''0.0000307'
$ws.Range('E12').Value = '  -3.30%  '
$ws.Range('D13').Value = '''9.53'
$ws.Range('E13').Value = '  +0.60%  '
$ws.Range('D14').Value = '''4.059.32'
$ws.Range('E14').Value = '  -1.42%  '
$ws.Range('D15').Value = '''607.89'
$ws.Range('E15').Value = '  +5.04%  '
$ws.Range('D16').Value = '''69.681.61'
$ws.Range('E16').Value = '  -1.52%  '
$ws.Range('D17').Value = '''18.96'
$ws.Range('E17').Value = '  -0.88%  '
$ws.Range('D18').Value = '''12.56'
$ws.Range('E18').Value = '  -1.51%  '
$ws.Range('D19').Value = '''3.506.49'
$ws.Range('E19').Value = '  -1.56%  '
$ws.Range('E20').Value = '  -0.26%  '
$ws.Range('E21').Value = '  -1.69%  '
$ws.Range('D22').Value = '''17.15'
$ws.Range('E22').Value = '  -2.68%  '
$ws.Range('D23').Value = '''106.06'
$ws.Range('E23').Value = '  +12.30%  '
$ws.Range('E24').Value = '  +3.38%  '
$ws.Range('D25').Value = '''5.09'
$ws.Range('E25').Value = '  +3.96%  '
$ws.Range('D26').Value = '''3.06'
$ws.Range('E26').Value = '  +4.41%  '
$ws.Range('D27').Value = '''10.99'
$ws.Range('E27').Value = '  -1.59%  '
$ws.Range('D28').Value = '''9.67'
$ws.Range('E28').Value = '  +3.96%  '
$ws.Range('D29').Value = '''33.70'
$ws.Range('E29').Value = '  +3.21%  '
$ws.Range('D30').Value = '''6.98'
$ws.Range('D31').Value = '''12.56'
$ws.Range('E31').Value = '  +2.19%  '
$ws.Range('D32').Value = '''3.92'
$ws.Range('E32').Value = '  +3.61%  '
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('D34').Value = '''63.48'
$ws.Range('E34').Value = '  -0.11%  '
$ws.Range('D35').Value = '''3.15'
$ws.Range('E35').Value = '  -4.71%  '
$ws.Range('E36').Value = '  -0.11%  '
$ws.Range('E37').Value = '  +6.96%  '
$ws.Range('D38').Value = '''3.653.34'
$ws.Range('E38').Value = '  +0.60%  '
$ws.Range('E39').Value = '  -3.99%  '
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').Value = '''508.66'
$ws.Range('E40').Value = '  -4.12%  '
$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D41').Value = '''36.72'
$ws.Range('E41').Value = '  -4.36%  '
$ws.Range('D42').Value = '''0.0₃0779'
$ws.Range('E42').Value = '  -3.16%  '
$ws.Range('D43').Value = '''0.137'
$ws.Range('E43').Value = '  -1.09%  '
$ws.Range('D44').Value = '''0.0464'
$ws.Range('E44').Value = '  -0.27%  '
$ws.Range('D45').Value = '''2.90'
$ws.Range('E45').Value = '  -0.62%  '
$ws.Range('E46').Value = '  +2.63%  '
$ws.Range('D47').Value = '''3.34'
$ws.Range('E47').Value = '  -3.66%  '
$ws.Range('E48').Value = '  +0.45%  '
$ws.Range('D49').Value = '''8.74'
$ws.Range('E49').Value = '  -5.83%  '
$ws.Range('D50').Value = '''131.86'
$ws.Range('E50').Value = '  -3.18%  '
$ws.Range('B51').Value = 'OceanProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean'
$ws.Range('D51').Value = '''1.35'
$ws.Range('E51').Value = '  -6.72%  '
